# Generate Report for Handback
#
# - Overview sheet: status moves from "Ready for handoff" to
#   "Handed back: in sync with en-US" for both language rows, and the two
#   status columns are widened to fit the new text.
# - zh-cn / de-de detail sheets: the "Latest Target File" (I) and
#   "Latest Handback File" (J) columns get populated now that handback
#   happened (previously blank), with I also becoming a hyperlink to the
#   source .md file (mirroring column A's hyperlink). The de-de sheet also
#   records the "Latest Handback DateTime" (K).
# - Column widths for the newly-populated columns are widened accordingly.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e65a1aae23c6bb1716c209cec7d1459c662d77a6/e2e/"
$mdFile1 = "1a1abcc3-fd4b-465c-be55-9db8d02ad460.md"
$mdFile2 = "f5cbce62-8268-45c0-82ff-cff954740dcc.md"

$handbackStatus = "Handed back: in sync with en-US"
$handbackDateTime = "2016-11-08 23:26:54"

# ---------------------------------------------------------------------
# Overview sheet: update status text + widen the zh-cn/de-de columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $handbackStatus
$wsOverview.Range("F2").Value = $handbackStatus
$wsOverview.Range("E3").Value = $handbackStatus
$wsOverview.Range("F3").Value = $handbackStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: fill in Latest Target File / Latest Handback File
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("J2").Value = "1a1abcc3-fd4b-465c-be55-9db8d02ad460.288a10e3611cedfb0d7d17042f5b7b164b243c61.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "f5cbce62-8268-45c0-82ff-cff954740dcc.29da179b020907edbda67fd99627d565b58f589e.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), ($repoBase + $mdFile1), "", "", $mdFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), ($repoBase + $mdFile2), "", "", $mdFile2)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: fill in Latest Target File / Latest Handback File /
# Latest Handback DateTime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("J2").Value = "1a1abcc3-fd4b-465c-be55-9db8d02ad460.288a10e3611cedfb0d7d17042f5b7b164b243c61.de-de.xlf"
$wsDeDe.Range("K2").Value = $handbackDateTime
$wsDeDe.Range("J3").Value = "f5cbce62-8268-45c0-82ff-cff954740dcc.29da179b020907edbda67fd99627d565b58f589e.de-de.xlf"
$wsDeDe.Range("K3").Value = $handbackDateTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), ($repoBase + $mdFile1), "", "", $mdFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), ($repoBase + $mdFile2), "", "", $mdFile2)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated: status + target/handback files + datetime updated."
